# "Revised Data for consistency"
# The "motivation" category previously labelled "Preventative Health"
# is renamed to "Preventative" everywhere it occurs on the sheet
# (rows 3, 8 and 13 — one row per location: Europe, Australia, United States).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Preventative"
$ws.Range("B8").Value = "Preventative"
$ws.Range("B13").Value = "Preventative"
